# wpiintlstudentsmap2 update:
#  - Row 71 (Huang, Chen S.) was missing its city/country ("sname"/"country")
#    values even though the "opcountry" column (F) already said "China".
#    Fill in the blanks to match.
#  - Row 185 (a near-empty stub row for "Ynan, Z." with no data beyond the
#    name and a blank country) is removed entirely, shifting all the
#    following rows up by one.
#  - Turn on AutoFilter over the data range and restore the last-used
#    selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B71").Value = "China"
$ws.Range("C71").Value = "China"

$ws.Rows(185).Delete()

$ws.Range("D24").Select()
$ws.Range("A1:I189").AutoFilter()
